$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B72 to be a numeric value instead of text
$ws.Range("B72").Value = 3

# Add new row 73 with data
$ws.Range("A73").Value = "Ying Tang"

# B73 must stay a text value ("4"), not get coerced to a number.
# Mark the cell as Text before assigning, then reset the style back to
# Normal so no stray number-format style is left referenced on the cell.
$ws.Range("B73").NumberFormat = "@"
$ws.Range("B73").Value = "4"
$ws.Range("B73").Style = "Normal"

$ws.Range("C73").Value = "We can of course remove"
$ws.Range("D73").Value = "ACK"
$ws.Range("E73").Value = "WRI"
$ws.Range("F73").Value = "a447d1c6-38e7-4648-9ebb-727dbeed5375"
$ws.Range("G73").Value = "SkFAWax0-_annotated.xlsx"
$ws.Range("H73").Value = "We can of course remove this part without taking away nothing from the paper's clarity, technical novelty and experimental success."
